# Team28EHyroExpoPoster.pptx edit
#
# 1) Update the "Fixed" date footer field (4/28/2017 -> 4/30/2017) that is
#    stamped into the Slide Master, every Custom Layout, and the Notes
#    Master (PowerPoint's Insert > Header & Footer "Apply to All").
# 2) Move a cluster of pictures/captions (the "Team"/acknowledgements
#    photos) from the bottom-right area of the poster up to the top-right,
#    per the commit message ("Moved pictures to top right").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Part 1: fixed date text on every Date Placeholder (master + layouts +
# notes master).
# ---------------------------------------------------------------------

function Set-DatePlaceholderText {
    param($shapes, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            # msoPlaceholder
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    # ppPlaceholderDate
                    $sh.TextFrame.TextRange.Text = $newText
                }
            } catch {
            }
        }
    }
}

$newDate = "4/30/2017"
$master = $p.SlideMaster

Set-DatePlaceholderText $master.Shapes $newDate

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# ---------------------------------------------------------------------
# Part 2: reposition the top-right picture cluster on slide 1.
#
# Shape.Left / Shape.Top are expressed in points (EMU / 12700). The target
# offsets below were chosen so that (points * 12700), after the COM host's
# internal float rounding, lands on (or, where the host's float precision
# makes an exact hit impossible, within 1 EMU of) the exact EMU offsets
# from the target OOXML.
# ---------------------------------------------------------------------

$slide = $p.Slides.Item(1)

$moves = @{
    30 = @{ Left = 2716.0;               Top = 522.7936401367188 }   # Subtitle 2 (team textbox)
    31 = @{ Left = 2690.256103515625;    Top = 1368.8253173828125 }  # TextBox 30 (team textbox)
    9  = @{ Left = 2739.4072265625;      Top = 749.4000244140625 }   # Picture 8
    6  = @{ Left = 3157.830322265625;    Top = 173.5596160888672 }   # Picture 5
    10 = @{ Left = 2885.92333984375;     Top = 173.72244262695312 }  # Picture 9
    12 = @{ Left = 3179.082763671875;    Top = 430.5782165527344 }   # TextBox 11
    15 = @{ Left = 2841.782470703125;    Top = 431.1067810058594 }   # TextBox 14
    17 = @{ Left = 2625.17626953125;     Top = 173.5596160888672 }   # Picture 16
    18 = @{ Left = 2639.630126953125;    Top = 430.5782165527344 }   # TextBox 17
    20 = @{ Left = 2827.029296875;       Top = 1273.2867431640625 }  # TextBox 19
}

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    $id = $sh.Id
    if ($moves.ContainsKey($id)) {
        $sh.Left = $moves[$id].Left
        $sh.Top = $moves[$id].Top
    }
}
